$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2198.3333
$ws.Range("I43").Value = 1950
$ws.Range("J43").Value = 2322.5
$ws.Range("K43").Value = 1950
$ws.Range("L43").Value = 2322.5
$ws.Range("M43").Value = -1881
$ws.Range("N43").Value = -2460.5
$ws.Range("H64").Value = 3898.5227
$ws.Range("I64").Value = 3733.2917
$ws.Range("J64").Value = 4096.8
$ws.Range("K64").Value = 3733.2917
$ws.Range("L64").Value = 4096.8
$ws.Range("M64").Value = -3485.2917
$ws.Range("N64").Value = -4592.8
$ws.Range("H67").Value = 3898.5227
$ws.Range("I67").Value = 3733.2917
$ws.Range("J67").Value = 4096.8
$ws.Range("K67").Value = 3733.2917
$ws.Range("L67").Value = 4096.8
$ws.Range("M67").Value = -2875.2917
$ws.Range("N67").Value = -5812.8
$ws.Range("H87").Value = 40666.668
$ws.Range("J87").Value = 40666.668
$ws.Range("L87").Value = 40666.668
$ws.Range("N87").Value = -43162.668
$ws.Range("H90").Value = 40666.668
$ws.Range("J90").Value = 40666.668
$ws.Range("L90").Value = 122000.004
$ws.Range("N90").Value = -134480.004
$ws.Range("H116").Value = 3661.5
$ws.Range("I116").Value = 3860
$ws.Range("K116").Value = 3860
$ws.Range("M116").Value = -418

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 18058.2
$ws.Range("I28").Value = 10723.286
$ws.Range("K28").Value = 10723.286
$ws.Range("M28").Value = -10531.286
$ws.Range("H99").Value = 18058.2
$ws.Range("I99").Value = 10723.286
$ws.Range("K99").Value = 10723.286
$ws.Range("M99").Value = -7728.286
$ws.Range("H132").Value = 52582.145
$ws.Range("I132").Value = 96731.82000000001
$ws.Range("J132").Value = 4017.5
$ws.Range("K132").Value = 290195.46
$ws.Range("L132").Value = 12052.5
$ws.Range("M132").Value = -287665.46
$ws.Range("N132").Value = -17112.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 944.1111
$ws.Range("I16").Value = 919.8
$ws.Range("J16").Value = 974.5
$ws.Range("K16").Value = 919.8
$ws.Range("L16").Value = 974.5
$ws.Range("M16").Value = -632.8
$ws.Range("N16").Value = -1548.5
$ws.Range("H22").Value = 384.2857
$ws.Range("I22").Value = 136.66667
$ws.Range("J22").Value = 570
$ws.Range("K22").Value = 136.66667
$ws.Range("L22").Value = 570
$ws.Range("M22").Value = 213.33333
$ws.Range("N22").Value = -1270
$ws.Range("H31").Value = 2983.303
$ws.Range("I31").Value = 1343.3158
$ws.Range("J31").Value = 5209
$ws.Range("K31").Value = 1343.3158
$ws.Range("L31").Value = 5209
$ws.Range("M31").Value = -1048.3158
$ws.Range("N31").Value = -5799
$ws.Range("H34").Value = 2983.303
$ws.Range("I34").Value = 1343.3158
$ws.Range("J34").Value = 5209
$ws.Range("K34").Value = 1343.3158
$ws.Range("L34").Value = 5209
$ws.Range("M34").Value = -1141.3158
$ws.Range("N34").Value = -5613
$ws.Range("H99").Value = 60841.35
$ws.Range("I99").Value = 68361
$ws.Range("K99").Value = 68361
$ws.Range("M99").Value = -66863
$ws.Range("H105").Value = 904.9524
$ws.Range("I105").Value = 869.3889
$ws.Range("J105").Value = 1118.3334
$ws.Range("K105").Value = 869.3889
$ws.Range("L105").Value = 1118.3334
$ws.Range("M105").Value = 877.6111
$ws.Range("N105").Value = -4612.3334
$ws.Range("H113").Value = 944.1111
$ws.Range("I113").Value = 919.8
$ws.Range("J113").Value = 974.5
$ws.Range("K113").Value = 919.8
$ws.Range("L113").Value = 974.5
$ws.Range("M113").Value = 1250.2
$ws.Range("N113").Value = -5314.5
$ws.Range("H126").Value = 60841.35
$ws.Range("I126").Value = 68361
$ws.Range("K126").Value = 205083
$ws.Range("M126").Value = -202613

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1089.7949
$ws.Range("I5").Value = 437.4762
$ws.Range("J5").Value = 1850.8334
$ws.Range("K5").Value = 1312.4286
$ws.Range("L5").Value = 5552.5002
$ws.Range("M5").Value = -1200.4286
$ws.Range("N5").Value = -5776.5002
$ws.Range("H38").Value = 49.695652
$ws.Range("I38").Value = 40
$ws.Range("J38").Value = 71.85714
$ws.Range("K38").Value = 120
$ws.Range("L38").Value = 215.57142
$ws.Range("M38").Value = 227
$ws.Range("N38").Value = -909.57142
$ws.Range("H80").Value = 916.5454999999999
$ws.Range("J80").Value = 997.1429000000001
$ws.Range("L80").Value = 2991.4287
$ws.Range("N80").Value = -4863.4287
$ws.Range("H83").Value = 916.5454999999999
$ws.Range("J83").Value = 997.1429000000001
$ws.Range("L83").Value = 8974.286100000001
$ws.Range("N83").Value = -18334.2861
$ws.Range("H113").Value = 614.14813
$ws.Range("I113").Value = 617.35297
$ws.Range("J113").Value = 608.7
$ws.Range("K113").Value = 1852.05891
$ws.Range("L113").Value = 1826.1
$ws.Range("M113").Value = 317.9410899999998
$ws.Range("N113").Value = -6166.1
$ws.Range("H122").Value = 508.76923
$ws.Range("I122").Value = 399.04544
$ws.Range("J122").Value = 1112.25
$ws.Range("K122").Value = 3591.40896
$ws.Range("L122").Value = 10010.25
$ws.Range("M122").Value = -1141.40896
$ws.Range("N122").Value = -14910.25
$ws.Range("H135").Value = 1089.7949
$ws.Range("I135").Value = 437.4762
$ws.Range("J135").Value = 1850.8334
$ws.Range("K135").Value = 3937.2858
$ws.Range("L135").Value = 16657.5006
$ws.Range("M135").Value = -1402.2858
$ws.Range("N135").Value = -21727.5006

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 37150
$ws.Range("J64").Value = 37150
$ws.Range("L64").Value = 37150
$ws.Range("N64").Value = -37600
$ws.Range("H67").Value = 37150
$ws.Range("J67").Value = 37150
$ws.Range("L67").Value = 37150
$ws.Range("N67").Value = -38710

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 23239.6
$ws.Range("J63").Value = 23239.6
$ws.Range("L63").Value = 23239.6
$ws.Range("N63").Value = -24487.6
$ws.Range("H66").Value = 23239.6
$ws.Range("J66").Value = 23239.6
$ws.Range("L66").Value = 69718.79999999999
$ws.Range("N66").Value = -75958.79999999999
